$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add a new worksheet "Sheet2" as the last (3rd) sheet in the book
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet2"

# Column A is wide in the new sheet
$newSheet.Columns.Item(1).ColumnWidth = 67.17

# ------------------------------------------------------------------
# 2. Fill in the content of the new sheet.
#    The order in which cells receive their Value matters because it
#    controls the order new strings are appended to the shared string
#    table, so we set them in the same order the original author did.
# ------------------------------------------------------------------
$newSheet.Range("A1").Value = "דברים שחקרנו/פתרנו"
$newSheet.Range("A2").Value = "פרוייקט אב"
$newSheet.Range("A3").Value = "הערות על משימות ופרוייקטים"
$newSheet.Range("A4").Value = "חברת אם"
$newSheet.Range("A5").Value = "יבוא משימות"
$newSheet.Range("A6").Value = "דף רישום"
$newSheet.Range("A7").Value = "טיוטות של פרוייקטים"
$newSheet.Range("A8").Value = "התממשקות עם מערכת פיננסית"
$newSheet.Range("A12").Value = "דברים שאנחנו מצפים מאלון ואלה"
$newSheet.Range("A9").Value = "ספקים לפעילות - כרגע לא מצאנו דרך מובנית - נדרש אפיון של מה התוצאה הנדרשת"
$newSheet.Range("A13").Value = "איזה תהליכים מול המערכת הפיננסית רוצים לממשק? אפיון"
$newSheet.Range("A14").Value = "פידבק כללי - והחלטה כללית אם התכולה הקיימת מספיקה להתחיל הטמעה - ואם לא, אז מה חסר"
$newSheet.Range("A15").Value = "נדרש איפיון+תכנון של יבוא הנתונים מתוך האקסלים."
$newSheet.Range("A10").Value = "pdf font/alignment"

# ------------------------------------------------------------------
# 3. Formatting
# ------------------------------------------------------------------
$newSheet.Range("A1").Font.Bold = $true
$newSheet.Range("A12").Font.Bold = $true
$newSheet.Range("A10").HorizontalAlignment = -4152   # xlRight

# ------------------------------------------------------------------
# 4. View settings for the new sheet: right-to-left, this sheet is the
#    active / selected tab, selection on A10.
# ------------------------------------------------------------------
$newSheet.Activate()
$excel.ActiveWindow.DisplayRightToLeft = $true
$newSheet.Range("A10").Select()

# ------------------------------------------------------------------
# 5. Update the view of the "אפיון " sheet: it is no longer the
#    selected tab, and its selection moves to I6.
# ------------------------------------------------------------------
$specSheet = $wb.Worksheets.Item("אפיון ")
$specSheet.Activate()
$specSheet.Range("I6").Select()

# Re-activate the new sheet so it ends up as the active tab of the
# workbook (matching activeTab pointing at Sheet2).
$newSheet.Activate()
$newSheet.Range("A10").Select()
